$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170 (shifts the existing rows 170-280 down to 171-281)
$ws.Rows("170:170").Insert()

# Populate the newly inserted row 170 with the new record
$ws.Range("A170").Value = 5
$ws.Range("B170").Value = "Macroferia Regional de Talca"
$ws.Range("C170").Value = "Maule"
$ws.Range("D170").Value = 44603
$ws.Range("E170").Value = 7
$ws.Range("F170").Value = 100114013
$ws.Range("G170").Value = "Zanahoria"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 400
$ws.Range("K170").Value = 7500
$ws.Range("L170").Value = 7500
$ws.Range("M170").Value = 7500
$ws.Range("N170").Value = "$/saco 20 kilos"
$ws.Range("O170").Value = "Región de Ñuble"
$ws.Range("P170").Value = 375
$ws.Range("Q170").Value = 20
$ws.Range("R170").Value = "Hortaliza"
